$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-77: ticker/column changes ---
$cellUpdates = @(
    @("B2", "NSE:5PAISA"),
    @("C2", "NSE:IVP"),
    @("D2", "NSE:CANBK"),
    @("E2", "NSE:PATANJALI"),
    @("F2", "NSE:BANKINDIA"),
    @("B3", "NSE:ABFRL"),
    @("C3", "NSE:LINC"),
    @("D3", "NSE:INDIANB"),
    @("E3", "NSE:PNBHOUSING"),
    @("F3", "NSE:DLF"),
    @("B4", "NSE:ACCURACY"),
    @("C4", "NSE:MALUPAPER"),
    @("D4", "NSE:LODHA"),
    @("F4", "NSE:INDIANB"),
    @("B5", "NSE:ALKALI"),
    @("C5", "NSE:MINDTECK"),
    @("D5", $null),
    @("F5", "NSE:INOXWIND"),
    @("B6", "NSE:ALLCARGO"),
    @("C6", "NSE:PGEL"),
    @("D6", $null),
    @("F6", "NSE:OBEROIRLTY"),
    @("B7", "NSE:ALOKINDS"),
    @("C7", "NSE:PNBHOUSING"),
    @("F7", "NSE:PNB"),
    @("B8", "NSE:AMNPLST"),
    @("B9", "NSE:ASHOKA"),
    @("B10", "NSE:BAJAJHLDNG"),
    @("B11", "NSE:BANKETF"),
    @("B12", "NSE:BASF"),
    @("B13", "NSE:BBTC"),
    @("B14", "NSE:BEPL"),
    @("B15", "NSE:BHAGERIA"),
    @("B16", "NSE:BHARATGEAR"),
    @("B17", "NSE:BIOFILCHEM"),
    @("B18", "NSE:CENTRUM"),
    @("B19", "NSE:CIEINDIA"),
    @("B20", "NSE:COMPUSOFT"),
    @("B21", "NSE:CONFIPET"),
    @("B22", "NSE:CORALFINAC"),
    @("B23", "NSE:CREST"),
    @("B24", "NSE:DHUNINV"),
    @("B25", "NSE:DLF"),
    @("B26", "NSE:DPWIRES"),
    @("B27", "NSE:DTIL"),
    @("B28", "NSE:DVL"),
    @("B29", "NSE:DYNPRO"),
    @("B30", "NSE:EKC"),
    @("B31", "NSE:ELECTCAST"),
    @("B32", "NSE:ESTER"),
    @("B33", "NSE:FMGOETZE"),
    @("B34", "NSE:FOODSIN"),
    @("B35", "NSE:GOLDIAM"),
    @("B36", "NSE:GREENPOWER"),
    @("B37", "NSE:GRINFRA"),
    @("B38", "NSE:GROBTEA"),
    @("B39", "NSE:HAPPSTMNDS"),
    @("B40", "NSE:HDFCPVTBAN"),
    @("B41", "NSE:HDFCQUAL"),
    @("B42", "NSE:HEALTHY"),
    @("B43", "NSE:HINDCON"),
    @("B44", "NSE:HINDWAREAP"),
    @("B45", "NSE:HLEGLAS"),
    @("B46", "NSE:HPL"),
    @("B47", "NSE:HYBRIDFIN"),
    @("B48", "NSE:IGPL"),
    @("B49", "NSE:IMFA"),
    @("B50", "NSE:INDIANB"),
    @("B51", "NSE:INDOWIND"),
    @("B52", "NSE:INFIBEAM"),
    @("B53", "NSE:INOXWIND"),
    @("B54", "NSE:ITDC"),
    @("B55", "NSE:JAYBARMARU"),
    @("B56", "NSE:JAYSREETEA"),
    @("B57", "NSE:JCHAC"),
    @("B58", "NSE:JINDALPOLY"),
    @("B59", "NSE:KEEPLEARN"),
    @("B60", "NSE:KERNEX"),
    @("B61", "NSE:KHAICHEM"),
    @("B62", "NSE:KILITCH"),
    @("B63", "NSE:KNRCON"),
    @("B64", "NSE:KUANTUM"),
    @("B65", "NSE:LIKHITHA"),
    @("B66", "NSE:MAHSEAMLES"),
    @("B67", "NSE:MOTHERSON"),
    @("B68", "NSE:MUNJALSHOW"),
    @("B69", "NSE:NDTV"),
    @("B70", "NSE:NUCLEUS"),
    @("B71", "NSE:OBEROIRLTY"),
    @("B72", "NSE:OSWALAGRO"),
    @("B73", "NSE:PHARMABEES"),
    @("B74", "NSE:PIXTRANS"),
    @("B75", "NSE:PLAZACABLE"),
    @("B76", "NSE:PNB"),
    @("B77", "NSE:PNCINFRA"),
)

foreach ($u in $cellUpdates) {
    $addr = $u[0]
    $val = $u[1]
    if ($null -eq $val) {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value = $val
    }
}

# --- Append new rows 78-86, copying column-A formatting from row 77 ---
$styleSrc = $ws.Range("A77")
$newRows = @(
    @(78, 76, "NSE:POCL"),
    @(79, 77, "NSE:PPL"),
    @(80, 78, "NSE:PRITIKAUTO"),
    @(81, 79, "NSE:RADIOCITY"),
    @(82, 80, "NSE:RATNAMANI"),
    @(83, 81, "NSE:REPCOHOME"),
    @(84, 82, "NSE:RKDL"),
    @(85, 83, "NSE:RTNPOWER"),
    @(86, 84, "NSE:SAHYADRI"),
)

foreach ($nr in $newRows) {
    $rowNum = $nr[0]
    $aVal = $nr[1]
    $bVal = $nr[2]
    $aCell = $ws.Range("A" + $rowNum)
    $styleSrc.Copy($aCell)
    $aCell.Value = $aVal
    if ($null -ne $bVal) {
        $ws.Range("B" + $rowNum).Value = $bVal
    }
}
